$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - values updated
$ws.Range("B3").Value = 0.9881019837019136
$ws.Range("C3").Value = 0.985965364695828
$ws.Range("D3").Value = 0.7592113448373879

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, values updated
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9838804333730072
$ws.Range("C4").Value = 0.9795343078278971
$ws.Range("D4").Value = 0.718280879723101

# Row 5: AdaBoostRegressor -> MLPRegressor, values updated
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.845012341260969
$ws.Range("C5").Value = 0.8339604357175182
$ws.Range("D5").Value = 0.602278401327291
